$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("programs")
$ws.Columns("C:C").Delete() | Out-Null
$ws.Columns("C:C").Select() | Out-Null
